# Daily attendance processing - 2026-01-19 04:54:32
#
# The "Recorded By" column (G) lists the users/processes that touched each
# attendance record, comma-separated. Re-processing re-orders that list so
# the last two contributors swap places (the most recently-merged recorder
# now sorts immediately ahead of the one before it), while any earlier
# contributors keep their original position. Rows whose "Recorded By" has
# only a single contributor (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the rows whose "Recorded By" cell actually lists more than one
# contributor need to change - touch exactly those rows (and nothing else)
# so untouched / genuinely-empty cells elsewhere on the sheet are not
# disturbed by the edit.
$rowsToFix = @(
    2,3,4,5,6,7,8,10,11,12,13,14,15,17,18,19,20,21,22,26,
    28,29,30,31,32,33,34,36,37,38,39,40,41,43,44,45,46,47,48,50,52,
    54,55,56,57,58,59,60,62,63,64,65,66,67,69,70,71,72,73,74,76,78,
    80,81,82,83,84,85,86,87,90,92,93,94,96,99,101,
    106,107,108,109,110,111,112,113,116,118,119,120,122,125,127,
    132,133,134,135,136,137,138,139,142,144,145,146,148,151,153
)

foreach ($row in $rowsToFix) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    $parts = $text.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $n = $parts.Length
    $last = $parts[$n - 1]
    $secondLast = $parts[$n - 2]
    $parts[$n - 1] = $secondLast
    $parts[$n - 2] = $last

    $newText = [string]::Join(", ", $parts)
    $cell.Value = $newText
}
